$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 4 columns (A..D) and 4 rows of data:
#   A1 "Parent name" | B1 "Parent type" | C1 "Comment"
#   A2 "Default"      | B2 "List"        | C2 "Comment 01 on list" | D2 (empty, dated style)
#   A3 "Task 01"      | B3 "Task"        | C3 "Comment 02 on task 01"
#   A4 "Task 01"      | B4 "Task"        | C4 "Comment 03 on task 01" | D4 (empty, dated style)
#
# Target: a single column with just the comment text ("Comment", "Comment 01",
# "Comment 02", "Comment 03"), column width 42.86 (column C's current width),
# and every cell sharing the plain "s=1" cell style that column C already uses.
#
# Column C already carries style index 1 on every row, and already has the
# exact target width (42.86), so we overwrite its values in place (this keeps
# reusing the existing style - no new style entries get created) and then
# delete the columns to its left so it slides into column A, preserving its
# width and per-cell styling exactly.

$ws.Range("C1").Value = "Comment"
$ws.Range("C2").Value = "Comment 01"
$ws.Range("C3").Value = "Comment 02"
$ws.Range("C4").Value = "Comment 03"

# Drop the old A and B columns (Parent name / Parent type) - this shifts
# column C (our new data, width 42.86) into column A.
$ws.Columns("A:B").Delete()

# Drop the now left-over old D column (the empty dated cells D2/D4).
$ws.Columns("B:B").Delete()

Write-Host "Rebuilt sheet with Comment column only"
